# feat: add 2022-Q1 data
#
# The workbook has quarterly "fund holdings" sheets (2020-Q4 .. 2021-Q4)
# followed by a "总计" (totals) summary sheet. This change:
#   1. Turns the existing "总计" sheet into the new "2022-Q1" detail sheet
#      (same slot/position, populated with the 2022-Q1 fund holdings table).
#   2. Appends a brand-new "总计" sheet after it, containing the old totals
#      table plus a new leading row for 2022-Q1.

$wb = $excel.ActiveWorkbook

# A sheet that already has the "header / index column" direct formatting
# (bold, centered, thin-bordered) we want to reuse as a paste-format
# source, so every new cell ends up with byte-identical styling instead of
# the engine inventing a fresh (but visually-equivalent) style each time a
# brand-new cell is formatted from a blank slate.
$formatTemplate = $wb.Worksheets.Item("2021-Q4")

# ------------------------------------------------------------------
# 1. Rename the current totals sheet to "2022-Q1" and add a fresh
#    "总计" sheet right after it.
# ------------------------------------------------------------------
$totalsSheet = $wb.Worksheets.Item("总计")
$totalsSheet.Name = "2022-Q1"

$newTotals = $wb.Worksheets.Add($null, $totalsSheet)
$newTotals.Name = "总计"

# ------------------------------------------------------------------
# 2. Build the new "总计" sheet: same 4-column layout as before, with
#    a new row inserted at the top for 2022-Q1.
# ------------------------------------------------------------------
$newTotals.Range("B1").Value = "日期"
$newTotals.Range("C1").Value = "持有数量(只)"
$newTotals.Range("D1").Value = "持有市值(亿元)"

$totalsRows = @(
    @("2022-Q1", 21,  8.220000000000001),
    @("2021-Q4", 16,  4.48),
    @("2021-Q3", 19,  8.67),
    @("2021-Q2", 16,  7.76),
    @("2021-Q1", 10,  4.7),
    @("2020-Q4", 7,   4.78)
)

for ($i = 0; $i -lt $totalsRows.Count; $i++) {
    $r = $i + 2
    $row = $totalsRows[$i]
    $newTotals.Range("A$r").Value = $i
    $newTotals.Range("B$r").Value = $row[0]
    $newTotals.Range("C$r").Value = $row[1]
    $newTotals.Range("D$r").Value = $row[2]
}

# This sheet is brand new, so every cell starts from a blank slate -
# paste the header / index-column formatting in from the template sheet
# (same combination used throughout the workbook: bold, centered,
# thin border).
$formatTemplate.Range("B1:D1").Copy()
$newTotals.Range("B1:D1").PasteSpecial(-4122)

$formatTemplate.Range("A2").Copy()
$newTotals.Range("A2:A7").PasteSpecial(-4122)

$newTotals.Range("A1").Select()

# ------------------------------------------------------------------
# 3. Build the "2022-Q1" sheet: the fund-holdings detail table
#    (same 8-column layout used by the other quarterly sheets).
#    This reuses the sheet that used to be "总计", which already has
#    the B1:D1 / A2:A6 cells carrying the header/index formatting -
#    only the newly-introduced cells (E1:H1 and A7:A22) need the
#    format pasted in explicitly.
# ------------------------------------------------------------------
$q1 = $totalsSheet
$q1.Cells.ClearContents()

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# code, name, size, position, ratio, value, rank
$detailRows = @(
    @("516970", "广发中证基建工程交易型开放式指数证券投资基金", "59.01", "99.38", "6.08", "3.5878", 4),
    @("515900", "博时中证央企创新驱动ETF", "47.96", "99.42", "2.98", "1.4292", 7),
    @("165525", "信诚中证基建工程指数（LOF）", "17.06", "94.00", "5.73", "0.9775", 4),
    @("515680", "嘉实中证央企创新驱动ETF", "17.86", "99.22", "2.99", "0.5340", 7),
    @("515600", "广发中证央企创新驱动ETF", "17.26", "99.02", "2.99", "0.5161", 7),
    @("516950", "银华中证基建交易型开放式指数证券投资基金", "10.41", "97.55", "4.51", "0.4695", 5),
    @("159974", "富国中证央企创新驱动ETF", "5.37", "99.51", "2.99", "0.1606", 7),
    @("000029", "富国宏观策略灵活配置混合", "5.97", "90.27", "2.02", "0.1206", 10),
    @("169108", "东方红均衡优选两年定期开放混合", "7.16", "22.76", "1.53", "0.1095", 2),
    @("160639", "鹏华中证高铁产业指数（LOF）", "0.89", "94.72", "8.04", "0.0716", 5),
    @("510160", "南方中证南方小康产业ETF", "2.60", "99.43", "2.19", "0.0569", 10),
    @("009242", "中加核心智造混合A", "2.05", "65.71", "2.34", "0.0480", 10),
    @("011484", "申万菱信宜选混合A", "4.30", "22.45", "0.89", "0.0383", 5),
    @("673071", "西部利得新动力灵活配置混合A", "4.64", "27.46", "0.49", "0.0227", 10),
    @("673073", "西部利得新动力灵活配置混合C", "3.81", "27.46", "0.49", "0.0187", 10),
    @("167702", "德邦量化优选股票(LOF)A", "0.54", "83.48", "2.86", "0.0154", 2),
    @("007505", "华夏中证AH经济蓝筹股票指数A", "1.09", "93.35", "1.25", "0.0136", 8),
    @("167703", "德邦量化优选股票(LOF)C", "0.41", "83.48", "2.86", "0.0117", 2),
    @("011485", "申万菱信宜选混合C", "1.07", "22.45", "0.89", "0.0095", 5),
    @("007506", "华夏中证AH经济蓝筹股票指数C", "0.65", "93.35", "1.25", "0.0081", 8),
    @("009243", "中加核心智造混合C", "0.10", "65.71", "2.34", "0.0023", 10)
)

# Columns B-G hold text (fund codes / percentages stored as strings,
# matching the source data), so force Text format before writing them.
$q1.Range("B2:G22").NumberFormat = "@"

for ($i = 0; $i -lt $detailRows.Count; $i++) {
    $r = $i + 2
    $row = $detailRows[$i]
    $q1.Range("A$r").Value = $i
    $q1.Range("B$r").Value = $row[0]
    $q1.Range("C$r").Value = $row[1]
    $q1.Range("D$r").Value = $row[2]
    $q1.Range("E$r").Value = $row[3]
    $q1.Range("F$r").Value = $row[4]
    $q1.Range("G$r").Value = $row[5]
    $q1.Range("H$r").Value = $row[6]
}

# B1:D1 and A2:A6 already carry the right formatting (they're the same
# cells that held "日期"/"持有数量(只)"/"持有市值(亿元)" and the old
# index column); only the newly-extended cells need it pasted in.
$formatTemplate.Range("B1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

$formatTemplate.Range("A2").Copy()
$q1.Range("A7:A22").PasteSpecial(-4122)

$q1.Range("A1").Select()
